$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.908.41'
$ws.Range("E2").Value = '  +0.77%  '

$ws.Range("D3").Value = '1.631.95'

$ws.Range("E4").Value = '  +0.39%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.518'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("E7").Value = '  +0.45%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.73'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.21%  '

$ws.Range("E9").Value = '  +0.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0608'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0904'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.26%  '

$ws.Range("D12").Value = '1.867.65'
$ws.Range("E12").Value = '  +1.59%  '

$ws.Range("D13").Value = '1.628.99'
$ws.Range("E13").Value = '  +2.02%  '

$ws.Range("E14").Value = '  -0.21%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '9.31'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.00%  '

$ws.Range("D16").Value = '29.917.87'
$ws.Range("E16").Value = '  +0.77%  '

$ws.Range("E17").Value = '  +0.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.12'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.87%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '241.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.08%  '

$ws.Range("E20").Value = '  -0.33%  '

$ws.Range("E22").Value = '  +1.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.85%  '

$ws.Range("E24").Value = '  +2.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.40'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.03%  '

$ws.Range("E26").Value = '  -0.88%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.109'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.25%  '

$ws.Range("E29").Value = '  +0.27%  '

$ws.Range("E30").Value = '  +1.70%  '

$ws.Range("E31").Value = '  +3.10%  '

$ws.Range("E32").Value = '  +2.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.17'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.02%  '

$ws.Range("D34").Value = '1.422.43'
$ws.Range("E34").Value = '  -0.51%  '

$ws.Range("E35").Value = '  +3.79%  '

$ws.Range("E36").Value = '  -1.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.75'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.33%  '

$ws.Range("E38").Value = '  +0.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0170'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.07%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '75.47'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.74%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.555'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.54%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.98'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.04%  '

$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0499'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.828'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.41%  '

$ws.Range("E45").Value = '  +0.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.13%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '50.81'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.71%  '

$ws.Range("D48").Value = '1.775.94'
$ws.Range("E48").Value = '  +1.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.33'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '90.41'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.91%  '

$ws.Range("E51").Value = '  +14.31%  '
